$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet "展览": currently only has the header row. Add the 13 event rows
# that were scraped for this update (rows 2-14), matching sheet "全部类型".
# =====================================================================
$ws1 = $wb.Worksheets.Item("展览")

# Column A holds a bold/centered/bordered "index" style (same as the header).
# Copy that formatting down before filling in the index values so the new
# cells pick up the same style as the existing header cell.
$ws1.Cells.Item(1, 1).Copy()
$ws1.Range("A2:A14").PasteSpecial(-4122)

$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(2, 2).Value = "'2024-07-06"
$ws1.Cells.Item(2, 3).Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws1.Cells.Item(2, 4).Value = '亭洪路45号 百益上河城'
$ws1.Cells.Item(2, 5).Value = '2024.07.06 10:00-07.06 17:00'
$ws1.Cells.Item(2, 6).Value = 615
$ws1.Cells.Item(2, 7).Value = 50
$ws1.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws1.Cells.Item(2, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'

$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(3, 2).Value = "'2024-07-06"
$ws1.Cells.Item(3, 3).Value = '南宁·首届童话梦境Lolita茶会'
$ws1.Cells.Item(3, 4).Value = '明秀东路157号 利泰国际大酒店'
$ws1.Cells.Item(3, 5).Value = '2024.07.06 13:00-07.06 17:00'
$ws1.Cells.Item(3, 6).Value = 206
$ws1.Cells.Item(3, 7).Value = 88
$ws1.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85776'
$ws1.Cells.Item(3, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'

$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(4, 2).Value = "'2024-07-12"
$ws1.Cells.Item(4, 3).Value = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$ws1.Cells.Item(4, 4).Value = '民族大道106号 南宁国际会展中心'
$ws1.Cells.Item(4, 5).Value = '2024.07.12 09:30-07.14 17:00'
$ws1.Cells.Item(4, 6).Value = 560
$ws1.Cells.Item(4, 7).Value = 29.9
$ws1.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87182'
$ws1.Cells.Item(4, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/CsYbpZmU1719311879090.jpeg'

$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(5, 2).Value = "'2024-07-13"
$ws1.Cells.Item(5, 3).Value = '南宁·0713国乙ONLY'
$ws1.Cells.Item(5, 4).Value = '亭洪路45号 水明漾宴会中心'
$ws1.Cells.Item(5, 5).Value = '2024.07.13 09:30-07.13 21:00'
$ws1.Cells.Item(5, 6).Value = 527
$ws1.Cells.Item(5, 7).Value = 68
$ws1.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws1.Cells.Item(5, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'

$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(6, 2).Value = "'2024-07-14"
$ws1.Cells.Item(6, 3).Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws1.Cells.Item(6, 4).Value = '明秀东路157号 利泰国际大酒店'
$ws1.Cells.Item(6, 5).Value = '2024.07.14 09:00-07.14 18:00'
$ws1.Cells.Item(6, 6).Value = 293
$ws1.Cells.Item(6, 7).Value = 69
$ws1.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws1.Cells.Item(6, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'

$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(7, 2).Value = "'2024-07-20"
$ws1.Cells.Item(7, 3).Value = '南宁·AB动漫游戏嘉年华'
$ws1.Cells.Item(7, 4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws1.Cells.Item(7, 5).Value = '2024.07.20 09:30-07.21 17:00'
$ws1.Cells.Item(7, 6).Value = 2647
$ws1.Cells.Item(7, 7).Value = 60
$ws1.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws1.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'

$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(8, 2).Value = "'2024-07-20"
$ws1.Cells.Item(8, 3).Value = '横州·第二届海棠动漫游戏嘉年华'
$ws1.Cells.Item(8, 4).Value = '茉莉花大道 横州国际大酒店'
$ws1.Cells.Item(8, 5).Value = '2024.07.20 09:30-07.20 17:00'
$ws1.Cells.Item(8, 6).Value = 450
$ws1.Cells.Item(8, 7).Value = 30
$ws1.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws1.Cells.Item(8, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'

$ws1.Cells.Item(9, 1).Value = 8
$ws1.Cells.Item(9, 2).Value = "'2024-07-27"
$ws1.Cells.Item(9, 3).Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws1.Cells.Item(9, 4).Value = '民族大道106号 南宁国际会展中心'
$ws1.Cells.Item(9, 5).Value = '2024.07.27 09:30-07.28 17:30'
$ws1.Cells.Item(9, 6).Value = 7342
$ws1.Cells.Item(9, 7).Value = 55
$ws1.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws1.Cells.Item(9, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'

$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(10, 2).Value = "'2024-08-03"
$ws1.Cells.Item(10, 3).Value = '南宁·火影忍者only'
$ws1.Cells.Item(10, 4).Value = '厢竹大道65号 桔子酒店'
$ws1.Cells.Item(10, 5).Value = '2024.08.03 10:00-08.03 17:00'
$ws1.Cells.Item(10, 6).Value = 193
$ws1.Cells.Item(10, 7).Value = 68
$ws1.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws1.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'

$ws1.Cells.Item(11, 1).Value = 10
$ws1.Cells.Item(11, 2).Value = "'2024-08-03"
$ws1.Cells.Item(11, 3).Value = '南宁·蔚蓝档案only'
$ws1.Cells.Item(11, 4).Value = '亭洪路45号 百益上河城'
$ws1.Cells.Item(11, 5).Value = '2024.08.03 09:00-08.03 17:00'
$ws1.Cells.Item(11, 6).Value = 456
$ws1.Cells.Item(11, 7).Value = 68
$ws1.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws1.Cells.Item(11, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'

$ws1.Cells.Item(12, 1).Value = 11
$ws1.Cells.Item(12, 2).Value = "'2024-08-10"
$ws1.Cells.Item(12, 3).Value = '南宁·国乙only'
$ws1.Cells.Item(12, 4).Value = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$ws1.Cells.Item(12, 5).Value = '2024.08.10 10:00-08.10 17:00'
$ws1.Cells.Item(12, 6).Value = 18
$ws1.Cells.Item(12, 7).Value = 40
$ws1.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88227'
$ws1.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/3cFX9LLQ1719482186347.jpeg'

$ws1.Cells.Item(13, 1).Value = 12
$ws1.Cells.Item(13, 2).Value = "'2024-08-24"
$ws1.Cells.Item(13, 3).Value = '南宁·第二届北极光动漫展'
$ws1.Cells.Item(13, 4).Value = '民族大道106号 南宁国际会展中心'
$ws1.Cells.Item(13, 5).Value = '2024.08.24 09:00-08.25 17:00'
$ws1.Cells.Item(13, 6).Value = 191
$ws1.Cells.Item(13, 7).Value = 65
$ws1.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88276'
$ws1.Cells.Item(13, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg'

$ws1.Cells.Item(14, 1).Value = 13
$ws1.Cells.Item(14, 2).Value = "'2024-11-02"
$ws1.Cells.Item(14, 3).Value = '南宁·万圣漫控嘉年华10'
$ws1.Cells.Item(14, 4).Value = '亭洪路45号 百益上河城'
$ws1.Cells.Item(14, 5).Value = '2024.11.02 11:00-11.03 22:00'
$ws1.Cells.Item(14, 6).Value = 39
$ws1.Cells.Item(14, 7).Value = 50
$ws1.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws1.Cells.Item(14, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'

# =====================================================================
# Sheet "演出": the 4 existing events now have a non-zero "想去人数"
# (want-to-go count) reported by the platform.
# =====================================================================
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 13
$ws2.Cells.Item(3, 6).Value = 17
$ws2.Cells.Item(4, 6).Value = 2
$ws2.Cells.Item(5, 6).Value = 1

# =====================================================================
# Sheet "全部类型": same "想去人数" refresh for rows 2-17 (union of the
# 展览 + 演出 rows above), plus a data fix -- row 18 was an accidental
# duplicate of row 17 (南宁·第二届北极光动漫展); replace it with the
# event that used to live in row 19 (南宁·万圣漫控嘉年华10) and drop the
# now-redundant row 19.
# =====================================================================
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 615
$ws4.Cells.Item(3, 6).Value = 206
$ws4.Cells.Item(4, 6).Value = 560
$ws4.Cells.Item(5, 6).Value = 527
$ws4.Cells.Item(6, 6).Value = 293
$ws4.Cells.Item(7, 6).Value = 13
$ws4.Cells.Item(8, 6).Value = 17
$ws4.Cells.Item(9, 6).Value = 2647
$ws4.Cells.Item(10, 6).Value = 450
$ws4.Cells.Item(11, 6).Value = 7342
$ws4.Cells.Item(12, 6).Value = 193
$ws4.Cells.Item(13, 6).Value = 456
$ws4.Cells.Item(14, 6).Value = 18
$ws4.Cells.Item(15, 6).Value = 2
$ws4.Cells.Item(16, 6).Value = 1
$ws4.Cells.Item(17, 6).Value = 191

# Row 18: overwrite with the row-19 event, keep the existing index (A18).
$ws4.Cells.Item(18, 2).Value = "'2024-11-02"
$ws4.Cells.Item(18, 3).Value = '南宁·万圣漫控嘉年华10'
$ws4.Cells.Item(18, 4).Value = '亭洪路45号 百益上河城'
$ws4.Cells.Item(18, 5).Value = '2024.11.02 11:00-11.03 22:00'
$ws4.Cells.Item(18, 6).Value = 39
$ws4.Cells.Item(18, 7).Value = 50
$ws4.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws4.Cells.Item(18, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'

# Row 19 is now redundant (its data lives in row 18) -- remove it so the
# sheet dimension shrinks back to A1:I18.
$ws4.Rows.Item(19).Delete()

